# Generate Report for Handoff
#
# The handoff generator re-ran for the "3a01482a-..." source file set
# (rows 4-7 in both the zh-cn and de-de sheets correspond to files
# 3a01482a, 48be1755, d6200de6 and ff6e9a9d):
#   - Priority moved from "low" -> "ht" on every one of those rows
#     (zh-cn AND de-de sheets).
#   - The zh-cn "Latest Handoff Datetime" (column H) advanced from
#     2016-08-28 02:31:15 -> 2016-08-28 02:31:31 on those same rows.
#   - The shared "Latest HO Xliff Generate Date" / handoff timestamp for
#     that file (2016-08-28 02:31:20 -> 2016-08-28 02:31:35) appears in
#     the Overview sheet (column G, rows 4-7) and in the de-de sheet
#     (column H, rows 4-7) since it is the same underlying value.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: rows 4-7 share the "Latest HO Xliff Generate Date" value
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-08-28 02:31:35"
}

# zh-cn sheet: rows 4-7 -> Priority "low" -> "ht", Latest Handoff Datetime updated
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-28 02:31:31"
}

# de-de sheet: rows 4-7 -> Priority "low" -> "ht"; column H shares the
# same "2016-08-28 02:31:20" -> "2016-08-28 02:31:35" value as Overview!G
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-28 02:31:35"
}
